# Actualización automática 2025-07-03 11:34:50
# Update PRESUPUESTO (column G) values on the "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G2").Value = 2500
$ws.Range("G3").Value = 1000
$ws.Range("G4").Value = 750
$ws.Range("G11").Value = 2000
$ws.Range("G13").Value = 1000
$ws.Range("G14").Value = 2500
$ws.Range("G18").Value = 3000
$ws.Range("G19").Value = 750
$ws.Range("G21").Value = 1000
$ws.Range("G22").Value = 1500
$ws.Range("G26").Value = 4000
$ws.Range("G28").Value = 750
$ws.Range("G29").Value = 500
$ws.Range("G31").Value = 1000

# Row 32 holds the total of column G (rows 2-31), stored as a static value.
$ws.Range("G32").Value = 27050
